$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the hyperlinks so the two password cells carry a display value
# with the original plaintext password (mirrors the mailto: links already
# present in the workbook). Rewriting individual Hyperlink properties in
# place leaves stray duplicate entries, so the whole collection is cleared
# and re-added from scratch instead.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:mercedez@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:P@ssword123", "", "", "P@ssword123")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:mercede@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:P@ssword1234567890", "", "", "P@ssword1234567890")

# Re-adding hyperlinks re-applies hyperlink formatting with a freshly
# minted style, so restore the original "Hyperlink" cell style on the
# affected cells to keep their styling untouched.
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("D3").Style = "Hyperlink"

# Encrypt (base64) the password values shown in the cells, keeping the
# original plaintext only as the hyperlink's display text (set above).
$ws.Range("D2").Value = "UEBzc3dvcmQxMjM="
$ws.Range("D3").Value = "UEBzc3dvcmQxMjM0NTY3ODkw"

# Update the selected cell to match the saved view state.
$ws.Range("E5").Select()
